$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (rows 2-5) replacing the previous 6 data rows (2-7).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics (updated TPM-derived values).
$data = @(
    @("ECs","Fgf9","Fgfr4","FAPs",3,1,4.077011333333332,12.231034,0.9715624748044627,0.9715624748044628,2,0.6666666666666666,0.1817723333333333,0.5453170000000001,0.008966262009224884,0.008966262009224884,0.7410878630864443,6.669790767777999,0.008711283707427762,0.008711283707427762),
    @("ECs","Fgf9","Fgfr4","MuSCs",3,1,4.077011333333332,12.231034,0.9715624748044627,0.9715624748044628,3,1,20.09115,60.27345,0.9910337379907751,0.9910337379907752,81.91184624969998,737.2066162472998,0.9628511910970349,0.9628511910970351),
    @("FAPs","Fgf9","Fgfr4","FAPs",1,0.3333333333333333,0.1193336666666667,0.358001,0.02843752519553723,0.02843752519553723,2,0.6666666666666666,0.1817723333333333,0.5453170000000001,0.008966262009224884,0.008966262009224884,0.02169155903522222,0.195224031317,0.0002549783017971209,0.0002549783017971209),
    @("FAPs","Fgf9","Fgfr4","MuSCs",1,0.3333333333333333,0.1193336666666667,0.358001,0.02843752519553723,0.02843752519553723,3,1,20.09115,60.27345,0.9910337379907751,0.9910337379907752,2.39755059705,21.57795537345,0.02818254689374011,0.02818254689374012)
)

# Remove the two rows that are being dropped (old rows 6 and 7), shrinking the table to 4 data rows.
$ws.Rows.Item(7).Delete() | Out-Null
$ws.Rows.Item(6).Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
